$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.072.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").Value = "'3.124.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.69%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'585.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").Value = "'134.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.54%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'3.119.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.82%  "
$ws.Range("E10").Value = "  -6.37%  "
$ws.Range("E11").Value = "  -4.90%  "
$ws.Range("D12").Value = "'0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.29%  "
$ws.Range("D13").Value = "'0.0000231"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.90%  "
$ws.Range("D14").Value = "'33.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "'3.638.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.80%  "
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "'3.133.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.41%  "
$ws.Range("D18").Value = "'62.078.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("D19").Value = "'6.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.35%  "
$ws.Range("D20").Value = "'450.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.09%  "
$ws.Range("D21").Value = "'13.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("D22").Value = "'0.697"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.93%  "
$ws.Range("D23").Value = "'7.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.87%  "
$ws.Range("D24").Value = "'13.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.93%  "
$ws.Range("D25").Value = "'82.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'2.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "'6.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.68%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.38%  "
$ws.Range("E31").Value = "  -8.42%  "
$ws.Range("D32").Value = "'26.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.19%  "
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("D34").Value = "'2.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.53%  "
$ws.Range("E35").Value = "  -8.21%  "
$ws.Range("E36").Value = "  -4.31%  "
$ws.Range("D37").Value = "'50.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.08%  "
$ws.Range("D38").Value = "'0.0₃0690"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.15%  "
$ws.Range("D39").Value = "'0.0382"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.57%  "
$ws.Range("D40").Value = "'2.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("D41").Value = "'7.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("D42").Value = "'393.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.25%  "
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").Value = "'2.744.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.65%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.247"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.72%  "
$ws.Range("D47").Value = "'2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.93%  "
$ws.Range("D48").Value = "'124.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("E49").Value = "  -5.24%  "
$ws.Range("D50").Value = "'33.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.85%  "
$ws.Range("E51").Value = "  -3.94%  "
